$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 33290.668
$ws.Range("J108").Value = 33290.668
$ws.Range("L108").Value = 33290.668
$ws.Range("N108").Value = -40970.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3668.825
$ws.Range("I2").Value = 3979.8057
$ws.Range("K2").Value = 3979.8057
$ws.Range("M2").Value = -3866.8057
$ws.Range("H28").Value = 7484
$ws.Range("I28").Value = 1921.3334
$ws.Range("K28").Value = 1921.3334
$ws.Range("M28").Value = -1729.3334
$ws.Range("H45").Value = 1800.2273
$ws.Range("I45").Value = 1579.55
$ws.Range("K45").Value = 1579.55
$ws.Range("M45").Value = -1202.55
$ws.Range("H61").Value = 2212.3948
$ws.Range("I61").Value = 1274.9615
$ws.Range("J61").Value = 4243.5
$ws.Range("K61").Value = 1274.9615
$ws.Range("L61").Value = 4243.5
$ws.Range("M61").Value = -1062.9615
$ws.Range("N61").Value = -4667.5
$ws.Range("H74").Value = 1408.1389
$ws.Range("I74").Value = 947.8261
$ws.Range("J74").Value = 2222.5386
$ws.Range("K74").Value = 947.8261
$ws.Range("L74").Value = 2222.5386
$ws.Range("M74").Value = -73.8261
$ws.Range("N74").Value = -3970.5386
$ws.Range("H77").Value = 1408.1389
$ws.Range("I77").Value = 947.8261
$ws.Range("J77").Value = 2222.5386
$ws.Range("K77").Value = 4739.1305
$ws.Range("L77").Value = 11112.693
$ws.Range("M77").Value = -371.1305000000002
$ws.Range("N77").Value = -19848.693
$ws.Range("H99").Value = 7484
$ws.Range("I99").Value = 1921.3334
$ws.Range("K99").Value = 1921.3334
$ws.Range("M99").Value = 1073.6666
$ws.Range("H102").Value = 68660
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 68660
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 68660
$ws.Range("N102").Value = -71904
$ws.Range("M102").ClearContents()
$ws.Range("H110").Value = 1926
$ws.Range("I110").Value = 1803.2354
$ws.Range("J110").Value = 4013
$ws.Range("K110").Value = 1803.2354
$ws.Range("L110").Value = 4013
$ws.Range("M110").Value = 241.7646
$ws.Range("N110").Value = -8103
$ws.Range("H116").Value = 3668.825
$ws.Range("I116").Value = 3979.8057
$ws.Range("K116").Value = 3979.8057
$ws.Range("M116").Value = -1685.8057
$ws.Range("H122").Value = 1980.5714
$ws.Range("I122").Value = 2446.2856
$ws.Range("J122").Value = 1514.8572
$ws.Range("K122").Value = 7338.8568
$ws.Range("L122").Value = 4544.571599999999
$ws.Range("M122").Value = -4888.8568
$ws.Range("N122").Value = -9444.571599999999
$ws.Range("H132").Value = 7043928.5
$ws.Range("I132").Value = 11364778
$ws.Range("J132").Value = 2544.1482
$ws.Range("K132").Value = 34094334
$ws.Range("L132").Value = 7632.444600000001
$ws.Range("M132").Value = -34091804
$ws.Range("N132").Value = -12692.4446
$ws.Range("H133").Value = 54516.168
$ws.Range("J133").Value = 54516.168
$ws.Range("L133").Value = 54516.168
$ws.Range("N133").Value = -59576.168
$ws.Range("H136").Value = 2212.3948
$ws.Range("I136").Value = 1274.9615
$ws.Range("J136").Value = 4243.5
$ws.Range("K136").Value = 3824.8845
$ws.Range("L136").Value = 12730.5
$ws.Range("M136").Value = -1274.8845
$ws.Range("N136").Value = -17830.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3668.825
$ws.Range("I3").Value = 3979.8057
$ws.Range("K3").Value = 3979.8057
$ws.Range("M3").Value = -3865.8057
$ws.Range("H105").Value = 2369.3845
$ws.Range("I105").Value = 2479.1
$ws.Range("K105").Value = 2479.1
$ws.Range("M105").Value = -732.0999999999999
$ws.Range("H107").Value = 1921.0938
$ws.Range("J107").Value = 1811.0769
$ws.Range("L107").Value = 1811.0769
$ws.Range("N107").Value = -5651.0769
$ws.Range("H134").Value = 3288.9265
$ws.Range("I134").Value = 2199.9033
$ws.Range("J134").Value = 4201.3516
$ws.Range("K134").Value = 6599.7099
$ws.Range("L134").Value = 12604.0548
$ws.Range("M134").Value = -4064.7099
$ws.Range("N134").Value = -17674.0548

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3935.1177
$ws.Range("I31").Value = 1499.7297
$ws.Range("J31").Value = 5812.396
$ws.Range("K31").Value = 1499.7297
$ws.Range("L31").Value = 5812.396
$ws.Range("M31").Value = -1204.7297
$ws.Range("N31").Value = -6402.396
$ws.Range("H34").Value = 3935.1177
$ws.Range("I34").Value = 1499.7297
$ws.Range("J34").Value = 5812.396
$ws.Range("K34").Value = 1499.7297
$ws.Range("L34").Value = 5812.396
$ws.Range("M34").Value = -1297.7297
$ws.Range("N34").Value = -6216.396
$ws.Range("H58").Value = 2969.2598
$ws.Range("I58").Value = 2977.2769
$ws.Range("K58").Value = 2977.2769
$ws.Range("M58").Value = -2774.2769
$ws.Range("H99").Value = 2292
$ws.Range("I99").Value = 2063.5
$ws.Range("J99").Value = 2383.4
$ws.Range("K99").Value = 2063.5
$ws.Range("L99").Value = 2383.4
$ws.Range("M99").Value = -565.5
$ws.Range("N99").Value = -5379.4
$ws.Range("H105").Value = 3048.4443
$ws.Range("J105").Value = 1888.875
$ws.Range("L105").Value = 1888.875
$ws.Range("N105").Value = -5382.875
$ws.Range("H126").Value = 2292
$ws.Range("I126").Value = 2063.5
$ws.Range("J126").Value = 2383.4
$ws.Range("K126").Value = 6190.5
$ws.Range("L126").Value = 7150.200000000001
$ws.Range("M126").Value = -3720.5
$ws.Range("N126").Value = -12090.2
$ws.Range("H132").Value = 29057.809
$ws.Range("I132").Value = 1266.0264
$ws.Range("K132").Value = 3798.0792
$ws.Range("M132").Value = -1268.0792
$ws.Range("H136").Value = 2969.2598
$ws.Range("I136").Value = 2977.2769
$ws.Range("K136").Value = 8931.830699999999
$ws.Range("M136").Value = -6381.830699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 4635
$ws.Range("J76").Value = 4780.4
$ws.Range("L76").Value = 14341.2
$ws.Range("N76").Value = -15107.2
$ws.Range("H79").Value = 4635
$ws.Range("J79").Value = 4780.4
$ws.Range("L79").Value = 14341.2
$ws.Range("N79").Value = -16993.2
$ws.Range("H94").Value = 3958.3333
$ws.Range("I94").Value = 1833.3334
$ws.Range("J94").Value = 4666.6665
$ws.Range("K94").Value = 5500.0002
$ws.Range("L94").Value = 13999.9995
$ws.Range("M94").Value = -4824.0002
$ws.Range("N94").Value = -15351.9995
$ws.Range("H100").Value = 5925
$ws.Range("J100").Value = 5925
$ws.Range("L100").Value = 17775
$ws.Range("N100").Value = -19397
$ws.Range("H106").Value = 870741.7
$ws.Range("J106").Value = 6112.5
$ws.Range("L106").Value = 18337.5
$ws.Range("N106").Value = -20229.5
$ws.Range("H112").Value = 3080347.8
$ws.Range("I112").Value = 100000950
$ws.Range("J112").Value = 3503.1746
$ws.Range("K112").Value = 300002850
$ws.Range("L112").Value = 10509.5238
$ws.Range("M112").Value = -300001742
$ws.Range("N112").Value = -12725.5238
$ws.Range("H122").Value = 2569.4717
$ws.Range("I122").Value = 622.4583
$ws.Range("J122").Value = 21260.8
$ws.Range("K122").Value = 5602.1247
$ws.Range("L122").Value = 191347.2
$ws.Range("M122").Value = -3152.1247
$ws.Range("N122").Value = -196247.2
$ws.Range("H125").Value = 105254.1
$ws.Range("I125").Value = 377632.5
$ws.Range("J125").Value = 6207.409
$ws.Range("K125").Value = 1132897.5
$ws.Range("L125").Value = 18622.227
$ws.Range("M125").Value = -1127977.5
$ws.Range("N125").Value = -28462.227
$ws.Range("H137").Value = 28578006
$ws.Range("I137").Value = 2292.3333
$ws.Range("J137").Value = 71441576
$ws.Range("K137").Value = 6876.999899999999
$ws.Range("L137").Value = 214324728
$ws.Range("M137").Value = -1776.999899999999
$ws.Range("N137").Value = -214334928

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 13567.777
$ws.Range("I113").Value = 18868.334
$ws.Range("K113").Value = 18868.334
$ws.Range("M113").Value = -16698.334
$ws.Range("H122").Value = 2424.7646
$ws.Range("I122").Value = 1986.2307
$ws.Range("J122").Value = 3850
$ws.Range("K122").Value = 5958.6921
$ws.Range("L122").Value = 11550
$ws.Range("M122").Value = -3508.6921
$ws.Range("N122").Value = -16450
$ws.Range("H132").Value = 3451.4055
$ws.Range("I132").Value = 3430.6155
$ws.Range("J132").Value = 3500.5454
$ws.Range("K132").Value = 10291.8465
$ws.Range("L132").Value = 10501.6362
$ws.Range("M132").Value = -7761.8465
$ws.Range("N132").Value = -15561.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 27998.666
$ws.Range("J99").Value = 29998.4
$ws.Range("L99").Value = 29998.4
$ws.Range("N99").Value = -35988.4
$ws.Range("H132").Value = 3928.647
$ws.Range("I132").Value = 3767.5
$ws.Range("K132").Value = 11302.5
$ws.Range("M132").Value = -8772.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1058856.1
$ws.Range("I122").Value = 1587923
$ws.Range("J122").Value = 722.2222
$ws.Range("K122").Value = 4763769
$ws.Range("L122").Value = 2166.6666
$ws.Range("M122").Value = -4761319
$ws.Range("N122").Value = -7066.6666
$ws.Range("H126").Value = 1838760.8
$ws.Range("I126").Value = 1961218.1
$ws.Range("K126").Value = 5883654.300000001
$ws.Range("M126").Value = -5881184.300000001
$ws.Range("H132").Value = 1374.9412
$ws.Range("I132").Value = 1055.8334
$ws.Range("K132").Value = 3167.5002
$ws.Range("M132").Value = -637.5001999999999
$ws.Range("H136").Value = 16670.027
$ws.Range("I136").Value = 28336.55
$ws.Range("J136").Value = 2528.7878
$ws.Range("K136").Value = 85009.64999999999
$ws.Range("L136").Value = 7586.3634
$ws.Range("M136").Value = -82459.64999999999
$ws.Range("N136").Value = -12686.3634
